$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the quantity multiplier in B18 (reduce order quantity from 10 to 1)
$ws.Range("B18").Value = 1

# Move selection to H18 (total cell) to reflect where user ended up after edit
$ws.Range("H18").Select()
